$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = 0.997864610145206
$ws.Range("D2").Value = 0.0000403400576448685
$ws.Range("E2").Value = 331.1848431006515
$ws.Range("F2").Value = 1447.507791857052
$ws.Range("G2").Value = 1116.3229487564
$ws.Range("H2").Value = 49150.05459653677
$ws.Range("I2").Value = 4103.130676205455
$ws.Range("J2").Value = 122.7899562924199
$ws.Range("K2").Value = 11.33966805013822
$ws.Range("L2").Value = 0.0000403400576448685
$ws.Range("M2").Value = 0.3831039430169797
$ws.Range("O2").Value = 0.00258
$ws.Range("P2").Value = 0.4052909952313641
$ws.Range("Q2").Value = 0.4016441217882361
$ws.Range("S2").Value = 0.01310853227799338
$ws.Range("T2").Value = 0.3142723187782764
$ws.Range("U2").Value = 0.9978173426149631
$ws.Range("V2").Value = 0.9979118821538743
$ws.Range("W2").Value = 134.1296243425581
$ws.Range("AC2").Value = 11.978

# Row 3 updates
$ws.Range("C3").Value = 0.9978292451472734
$ws.Range("D3").Value = 0.00003583090569095404
$ws.Range("E3").Value = 331.1965809318972
$ws.Range("F3").Value = 1447.503005802524
$ws.Range("G3").Value = 1116.306424870626
$ws.Range("H3").Value = 48194.59125163723
$ws.Range("I3").Value = 4256.947837843626
$ws.Range("J3").Value = 113.0789815084083
$ws.Range("K3").Value = 11.09199890916404
$ws.Range("L3").Value = 0.00003583090569095404
$ws.Range("M3").Value = 0.3857432425476348
$ws.Range("P3").Value = 0.399786534923933
$ws.Range("Q3").Value = 0.4074675873316524
$ws.Range("S3").Value = 0.01162185150579935
$ws.Range("T3").Value = 0.3016542752120749
$ws.Range("U3").Value = 0.9977870018936271
$ws.Range("V3").Value = 0.997871491977972
$ws.Range("W3").Value = 124.1709804175723
$ws.Range("AC3").Value = 12.08

# Row 4 updates
$ws.Range("C4").Value = 0.9978095733888616
$ws.Range("D4").Value = 0.00003906085247667861
$ws.Range("E4").Value = 331.2031104534619
$ws.Range("F4").Value = 1447.521962537725
$ws.Range("G4").Value = 1116.318852084263
$ws.Range("H4").Value = 46080.55311799021
$ws.Range("I4").Value = 3817.579597394361
$ws.Range("J4").Value = 109.3331012484569
$ws.Range("K4").Value = 11.0728536290639
$ws.Range("L4").Value = 0.00003906085247667861
$ws.Range("M4").Value = 0.3896362373277412
$ws.Range("P4").Value = 0.4037236656929212
$ws.Range("Q4").Value = 0.4066422061589337
$ws.Range("S4").Value = 0.01270681042663709
$ws.Range("T4").Value = 0.2643135192520635
$ws.Range("U4").Value = 0.9977639730800651
$ws.Range("V4").Value = 0.9978551778659448
$ws.Range("W4").Value = 120.4059548775208
$ws.Range("AC4").Value = 19.977

# Row 5 updates
$ws.Range("C5").Value = 0.9977819301559424
$ws.Range("D5").Value = 0.00004021169415392708
$ws.Range("E5").Value = 331.2122863309248
$ws.Range("F5").Value = 1447.539256706479
$ws.Range("G5").Value = 1116.326970375554
$ws.Range("H5").Value = 44440.93230787192
$ws.Range("I5").Value = 3778.687317392505
$ws.Range("J5").Value = 106.9642679490747
$ws.Range("K5").Value = 11.41406198750633
$ws.Range("L5").Value = 0.00004021169415392707
$ws.Range("M5").Value = 0.3914127086273933
$ws.Range("P5").Value = 0.4092324959144185
$ws.Range("Q5").Value = 0.4091641516886602
$ws.Range("S5").Value = 0.0130969040856272
$ws.Range("T5").Value = 0.2628163049414113
$ws.Range("U5").Value = 0.9977351874418423
$ws.Range("V5").Value = 0.9978286772499295
$ws.Range("W5").Value = 118.378329936581
$ws.Range("AC5").Value = 12.044

# Row 6 updates
$ws.Range("C6").Value = 0.9977783725275983
$ws.Range("D6").Value = 0.00003762356284133659
$ws.Range("E6").Value = 331.2134672847822
$ws.Range("F6").Value = 1447.541012035804
$ws.Range("G6").Value = 1116.327544751022
$ws.Range("H6").Value = 45288.66273637601
$ws.Range("I6").Value = 3904.517386088566
$ws.Range("J6").Value = 110.8980180547772
$ws.Range("K6").Value = 12.56714227142768
$ws.Range("L6").Value = 0.00003762356284133659
$ws.Range("M6").Value = 0.3933162854106986
$ws.Range("P6").Value = 0.4163746000991697
$ws.Range("Q6").Value = 0.4102546201380286
$ws.Range("S6").Value = 0.01220602945670827
$ws.Range("T6").Value = 0.2668567001166197
$ws.Range("U6").Value = 0.9977340425561128
$ws.Range("V6").Value = 0.9978227064384776
$ws.Range("W6").Value = 123.4651603262049
$ws.Range("AC6").Value = 19.988

# Row 7 updates
$ws.Range("C7").Value = 0.9977801156790581
$ws.Range("D7").Value = 0.00003554953358906369
$ws.Range("E7").Value = 331.2128886450298
$ws.Range("F7").Value = 1447.539632846651
$ws.Range("G7").Value = 1116.326744201622
$ws.Range("H7").Value = 46370.33402072885
$ws.Range("I7").Value = 4282.108479381145
$ws.Range("J7").Value = 116.1397256563911
$ws.Range("K7").Value = 12.37561484100051
$ws.Range("L7").Value = 0.00003554953358906369
$ws.Range("M7").Value = 0.3925918093317671
$ws.Range("P7").Value = 0.4170020136156424
$ws.Range("Q7").Value = 0.4013845083479939
$ws.Range("S7").Value = 0.01149056550647157
$ws.Range("T7").Value = 0.3407948763191531
$ws.Range("U7").Value = 0.9977377599688388
$ws.Range("V7").Value = 0.9978224749855779
$ws.Range("W7").Value = 128.5153404973916
$ws.Range("AC7").Value = 12.039

# Row 8 updates
$ws.Range("C8").Value = 0.997763519810411
$ws.Range("D8").Value = 0.00003996281559342204
$ws.Range("E8").Value = 331.2183977315869
$ws.Range("F8").Value = 1447.534558137279
$ws.Range("G8").Value = 1116.316160405692
$ws.Range("H8").Value = 44105.47704015052
$ws.Range("I8").Value = 3941.970379392729
$ws.Range("J8").Value = 102.6333952862651
$ws.Range("K8").Value = 11.52098883377944
$ws.Range("L8").Value = 0.00003996281559342204
$ws.Range("M8").Value = 0.3949885259756325
$ws.Range("P8").Value = 0.4058012811822291
$ws.Range("Q8").Value = 0.4044871219574797
$ws.Range("S8").Value = 0.01302433137593426
$ws.Range("T8").Value = 0.3286078745419824
$ws.Range("U8").Value = 0.9977171781464422
$ws.Range("V8").Value = 0.9978098657795067
$ws.Range("W8").Value = 114.1543841200446
$ws.Range("AC8").Value = 12.059

# Row 9 updates
$ws.Range("C9").Value = 0.9977587176510911
$ws.Range("D9").Value = 0.0000411619038602536
$ws.Range("E9").Value = 331.2199918680124
$ws.Range("F9").Value = 1447.516894809508
$ws.Range("G9").Value = 1116.296902941496
$ws.Range("H9").Value = 43261.36597323853
$ws.Range("I9").Value = 3858.442438382513
$ws.Range("J9").Value = 108.0931650600248
$ws.Range("K9").Value = 11.88368840996711
$ws.Range("L9").Value = 0.0000411619038602536
$ws.Range("M9").Value = 0.392970700582273
$ws.Range("P9").Value = 0.4140357776661653
$ws.Range("Q9").Value = 0.4052672146144209
$ws.Range("S9").Value = 0.01338922658632062
$ws.Range("T9").Value = 0.3383168108910171
$ws.Range("U9").Value = 0.9977106447972814
$ws.Range("V9").Value = 0.9978067951377283
$ws.Range("W9").Value = 119.9768534699919
$ws.Range("AC9").Value = 12.094

# Row 10 updates
$ws.Range("C10").Value = 0.9977601746701813
$ws.Range("D10").Value = 0.00004086339614387934
$ws.Range("E10").Value = 331.2195081908087
$ws.Range("F10").Value = 1447.52473729567
$ws.Range("G10").Value = 1116.305229104861
$ws.Range("H10").Value = 43156.96590535386
$ws.Range("I10").Value = 3811.220077453418
$ws.Range("J10").Value = 103.7944535080573
$ws.Range("K10").Value = 12.61453751755598
$ws.Range("L10").Value = 0.00004086339614387934
$ws.Range("M10").Value = 0.3925630102433947
$ws.Range("P10").Value = 0.4168813180695979
$ws.Range("Q10").Value = 0.4062179138014888
$ws.Range("S10").Value = 0.01331299985468358
$ws.Range("T10").Value = 0.3314984253568839
$ws.Range("U10").Value = 0.9977127228504987
$ws.Range("V10").Value = 0.9978076310037529
$ws.Range("W10").Value = 116.4089910256133
$ws.Range("AC10").Value = 0
